$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "The jar near the sticky lid"
$ws.Range("C3").Value = "The jar near the sticky cookie"
$ws.Range("C4").Value = "The jar near the stale bagel"
$ws.Range("D2").Value = "The jar near the sticky lids"
$ws.Range("D3").Value = "The jar near the sticky cookies"
$ws.Range("D4").Value = "The jar near the stale bagels"
